$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028382726690068
$ws.Range("D2").Value = 1.031956381641022
$ws.Range("E2").Value = 1.028104918743613
$ws.Range("F2").Value = 1.026932373238568
$ws.Range("I2").Value = 1.031040219857644
$ws.Range("J2").Value = 1.033535352532168
$ws.Range("K2").Value = 1.034763002514901
$ws.Range("L2").Value = 1.03092270516103
$ws.Range("M2").Value = 1.029753576994986
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030325156207344
$ws.Range("D3").Value = 1.033444566559109
$ws.Range("E3").Value = 1.029806281502341
$ws.Range("F3").Value = 1.029524038627695
$ws.Range("I3").Value = 1.031561218825675
$ws.Range("J3").Value = 1.035113283085584
$ws.Range("K3").Value = 1.036058040677981
$ws.Range("L3").Value = 1.032429523201192
$ws.Range("M3").Value = 1.032148041176021
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031575916713346
$ws.Range("D4").Value = 1.034402243388289
$ws.Range("E4").Value = 1.030901306094183
$ws.Range("F4").Value = 1.031193934884516
$ws.Range("I4").Value = 1.031894252084095
$ws.Range("J4").Value = 1.036128069497138
$ws.Range("K4").Value = 1.036890255375326
$ws.Range("L4").Value = 1.033398233619777
$ws.Range("M4").Value = 1.033690114697445
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032100302485043
$ws.Range("D5").Value = 1.03480360995001
$ws.Range("E5").Value = 1.031360276552743
$ws.Range("F5").Value = 1.031894310293911
$ws.Range("I5").Value = 1.032033290655612
$ws.Range("J5").Value = 1.036553218528711
$ws.Range("K5").Value = 1.037238760779641
$ws.Range("L5").Value = 1.033803996914007
$ws.Range("M5").Value = 1.034336697554263
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032188265985439
$ws.Range("D6").Value = 1.034870929061682
$ws.Range("E6").Value = 1.031437259709995
$ws.Range("F6").Value = 1.032011811068504
$ws.Range("I6").Value = 1.032056579333374
$ws.Range("J6").Value = 1.036624517657347
$ws.Range("K6").Value = 1.037297197382472
$ws.Range("L6").Value = 1.033872040170093
$ws.Range("M6").Value = 1.034445162956073
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031582929173871
$ws.Range("D7").Value = 1.034407611312612
$ws.Range("E7").Value = 1.030907444256973
$ws.Range("F7").Value = 1.031203299751472
$ws.Range("I7").Value = 1.03189611371675
$ws.Range("J7").Value = 1.036133756088438
$ws.Range("K7").Value = 1.036894917427127
$ws.Range("L7").Value = 1.033403661239287
$ws.Range("M7").Value = 1.033698761009374
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029040466392303
$ws.Range("D8").Value = 1.032460427482943
$ws.Range("E8").Value = 1.028681134605332
$ws.Range("F8").Value = 1.027809734143545
$ws.Range("I8").Value = 1.031217146423547
$ws.Range("J8").Value = 1.034069929831923
$ws.Range("K8").Value = 1.035201872853842
$ws.Range("L8").Value = 1.03143326196043
$ws.Range("M8").Value = 1.030564336566534
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024512010649539
$ws.Range("D9").Value = 1.02898777636902
$ws.Range("E9").Value = 1.024711875118401
$ws.Range("F9").Value = 1.021773441640437
$ws.Range("I9").Value = 1.029988953776024
$ws.Range("J9").Value = 1.030384198110756
$ws.Range("K9").Value = 1.032173399892303
$ws.Range("L9").Value = 1.027911721328044
$ws.Range("M9").Value = 1.024983137572472
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021458538040969
$ws.Range("D10").Value = 1.02664332777225
$ws.Range("E10").Value = 1.022032857138948
$ws.Range("F10").Value = 1.017708309824393
$ws.Range("I10").Value = 1.029148153623308
$ws.Range("J10").Value = 1.027892378560218
$ws.Range("K10").Value = 1.030122693147082
$ws.Range("L10").Value = 1.025529112329011
$ws.Range("M10").Value = 1.021220561254592
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020127711578314
$ws.Range("D11").Value = 1.025620862726974
$ws.Range("E11").Value = 1.020864622950147
$ws.Range("F11").Value = 1.015937672331922
$ws.Range("I11").Value = 1.028778711876489
$ws.Range("J11").Value = 1.026804783745611
$ws.Range("K11").Value = 1.02922687481471
$ws.Range("L11").Value = 1.024488757943796
$ws.Range("M11").Value = 1.019580779263869
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01963204376406
$ws.Range("D12").Value = 1.025239947503323
$ws.Range("E12").Value = 1.020429421781899
$ws.Range("F12").Value = 1.015278355770071
$ws.Range("I12").Value = 1.028640665187868
$ws.Range("J12").Value = 1.026399473740668
$ws.Range("K12").Value = 1.028892921344701
$ws.Range("L12").Value = 1.024100988773694
$ws.Range("M12").Value = 1.01897004916031
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019738427457443
$ws.Range("D13").Value = 1.025321706544415
$ws.Range("E13").Value = 1.02052283180804
$ws.Range("F13").Value = 1.015419855815476
$ws.Range("I13").Value = 1.028670313967416
$ws.Range("J13").Value = 1.026486474794945
$ws.Range("K13").Value = 1.028964610547932
$ws.Range("L13").Value = 1.024184227532332
$ws.Range("M13").Value = 1.019101128046595
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020086767081912
$ws.Range("D14").Value = 1.025589399298964
$ws.Range("E14").Value = 1.020828675135185
$ws.Range("F14").Value = 1.015883206473769
$ws.Range("I14").Value = 1.028767317688402
$ws.Range("J14").Value = 1.02677130800013
$ws.Range("K14").Value = 1.029199294892947
$ws.Range("L14").Value = 1.024456732253611
$ws.Range("M14").Value = 1.019530329944396
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020301211910651
$ws.Range("D15").Value = 1.025754183534658
$ws.Range("E15").Value = 1.021016946365465
$ws.Range("F15").Value = 1.016168475258108
$ws.Range("I15").Value = 1.028826975892427
$ws.Range("J15").Value = 1.026946626030518
$ws.Range("K15").Value = 1.029343730938973
$ws.Range("L15").Value = 1.024624453482798
$ws.Range("M15").Value = 1.019794556136921
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021546674941731
$ws.Range("D16").Value = 1.026711029079981
$ws.Range("E16").Value = 1.022110213241758
$ws.Range("F16").Value = 1.017825596880881
$ws.Range("I16").Value = 1.029172558031801
$ws.Range("J16").Value = 1.027964374157586
$ws.Range("K16").Value = 1.030181977962747
$ws.Range("L16").Value = 1.025597971791251
$ws.Range("M16").Value = 1.021329160968484
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022325576679331
$ws.Range("D17").Value = 1.027309257450311
$ws.Range("E17").Value = 1.022793770502555
$ws.Range("F17").Value = 1.018862237401202
$ws.Range("I17").Value = 1.029387885970067
$ws.Range("J17").Value = 1.028600449105458
$ws.Range("K17").Value = 1.030705667613679
$ws.Range("L17").Value = 1.026206290206441
$ws.Range("M17").Value = 1.022288911694161
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022779064742774
$ws.Range("D18").Value = 1.027657491255012
$ws.Range("E18").Value = 1.023191688920832
$ws.Range("F18").Value = 1.019465892055272
$ws.Range("I18").Value = 1.029512965892232
$ws.Range("J18").Value = 1.028970631418185
$ws.Range("K18").Value = 1.031010371515885
$ws.Range("L18").Value = 1.026560277864123
$ws.Range("M18").Value = 1.022847702800694
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022933552496067
$ws.Range("D19").Value = 1.027776111663344
$ws.Range("E19").Value = 1.023327236037993
$ws.Range("F19").Value = 1.019671554811429
$ws.Range("I19").Value = 1.029555527568042
$ws.Range("J19").Value = 1.029096714343712
$ws.Range("K19").Value = 1.031114140364525
$ws.Range("L19").Value = 1.026680837999662
$ws.Range("M19").Value = 1.023038065580387
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022242094251393
$ws.Range("D20").Value = 1.027245146097292
$ws.Range("E20").Value = 1.022720513133791
$ws.Range("F20").Value = 1.018751119566556
$ws.Range("I20").Value = 1.029364836910874
$ws.Range("J20").Value = 1.028532290318682
$ws.Range("K20").Value = 1.030649559003193
$ws.Range("L20").Value = 1.026141109916578
$ws.Range("M20").Value = 1.022186044933655
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019984227046583
$ws.Range("D21").Value = 1.025510601778107
$ws.Range("E21").Value = 1.020738647155386
$ws.Range("F21").Value = 1.015746806536324
$ws.Range("I21").Value = 1.028738775238941
$ws.Range("J21").Value = 1.026687468677798
$ws.Range("K21").Value = 1.029130219735514
$ws.Range("L21").Value = 1.024376523530688
$ws.Range("M21").Value = 1.019403986461453
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018556845553122
$ws.Range("D22").Value = 1.024413494438553
$ws.Range("E22").Value = 1.019485221330257
$ws.Range("F22").Value = 1.013848455577009
$ws.Range("I22").Value = 1.028340397064731
$ws.Range("J22").Value = 1.025519850498669
$ws.Range("K22").Value = 1.028167955501564
$ws.Range("L22").Value = 1.023259316095616
$ws.Range("M22").Value = 1.017645267233159
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019314277503414
$ws.Range("D23").Value = 1.024995720621448
$ws.Range("E23").Value = 1.020150394360884
$ws.Range("F23").Value = 1.014855720324384
$ws.Range("I23").Value = 1.028552039314348
$ws.Range("J23").Value = 1.026139568821822
$ws.Range("K23").Value = 1.028678742366845
$ws.Range("L23").Value = 1.023852313861153
$ws.Range("M23").Value = 1.018578519526368
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022279818936753
$ws.Range("D24").Value = 1.02727411742821
$ws.Range("E24").Value = 1.022753617421858
$ws.Range("F24").Value = 1.018801332020944
$ws.Range("I24").Value = 1.029375253386271
$ws.Range("J24").Value = 1.028563090902105
$ws.Range("K24").Value = 1.030674914402184
$ws.Range("L24").Value = 1.026170564660129
$ws.Range("M24").Value = 1.02223252912944
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025688655906901
$ws.Range("D25").Value = 1.029890603266122
$ws.Range("E25").Value = 1.025743681221349
$ws.Range("F25").Value = 1.023340958592643
$ws.Range("I25").Value = 1.030310301287518
$ws.Range("J25").Value = 1.031343030326792
$ws.Range("K25").Value = 1.032961820654317
$ws.Range("L25").Value = 1.028828153701717
$ws.Range("M25").Value = 1.026433164878064
